$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Move the "A 24922-2019" record (currently row 43) to the top of
#    the data block (row 4), pushing the rows that were 4..42 down
#    by one (they become 5..43). Its statistics / species / links
#    have also been refreshed as part of the move.
# -----------------------------------------------------------------

# Insert a new blank row at row 4 (formatting is inherited from row 3,
# which already uses the correct date / wrap-text styles).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the refreshed "A 24922-2019" data.
$ws.Range("A4").Value = "A 24922-2019"
$ws.Range("B4").Value = 43599
$ws.Range("C4").Value = 45178
$ws.Range("D4").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E4").Value = "VÄNNÄS"
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = "Garnlav`nJärpe`nTretåig hackspett`nUllticka"

$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/artfynd/A 24922-2019.xlsx")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/kartor/A 24922-2019.png")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/klagomål/A 24922-2019.docx")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/klagomålsmail/A 24922-2019.docx")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/tillsyn/A 24922-2019.docx")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/tillsynsmail/A 24922-2019.docx")'

# Keep the row height consistent with the rest of the sheet (Excel
# would otherwise auto-grow it because of the wrapped, multi-line
# species list).
$ws.Rows.Item(4).RowHeight = 15

# The old "A 24922-2019" row has shifted from row 43 to row 44 because
# of the insert above; remove that now-duplicate row so every other
# record keeps its original relative order.
$ws.Rows.Item(44).Delete()

# -----------------------------------------------------------------
# 2. Refresh the "Förändrad" (changed) date stamp in column C for
#    every data row from 45177 to 45178.
# -----------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value = 45178
}
